$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the "目的" description text - remove "システム開発"
$ws.Range("A4").Value = "管理者が本をサーバー上にアプロードし、ユーザー側でウェブ上本を買えるサービスに関する要件定義書になります。"

# 2) Rename "配送先住所画面" -> "配送先画面"
$ws.Range("A54").Value = "７．配送先画面"

# 3) Rename "請求先住所画面" -> "請求先画面"
$ws.Range("A55").Value = "８．請求先画面"

# 4) Update the sheet view: drop the frozen/scrolled topLeftCell, and move
#    the active selection to A40.
$ws.Activate()
$ws.Range("A40").Select()
